$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C2").Value = 4.916123640948215
$ws.Range("D2").Value = 8.969400099930523
$ws.Range("E2").Value = 13.58350496176993
$ws.Range("F2").Value = 33.7675643087357
$ws.Range("G2").Value = 3.654299887884765
$ws.Range("J2").Value = 9.948131869743699
$ws.Range("K2").Value = 15.95505938484501
$ws.Range("N2").Value = 17.41358868209829
$ws.Range("O2").Value = 25.44303966054887
$ws.Range("C3").Value = 4.750579149002285
$ws.Range("D3").Value = 8.923348493824665
$ws.Range("E3").Value = 13.52742842962071
$ws.Range("F3").Value = 33.77018483041511
$ws.Range("G3").Value = 3.656784836498327
$ws.Range("J3").Value = 9.955227870812518
$ws.Range("K3").Value = 15.50193656273498
$ws.Range("N3").Value = 17.45741776085186
$ws.Range("O3").Value = 25.49628641335583
$ws.Range("C4").Value = 4.647470269827523
$ws.Range("D4").Value = 8.896327848593273
$ws.Range("E4").Value = 13.49570476194992
$ws.Range("F4").Value = 33.78151763719489
$ws.Range("G4").Value = 3.65839104939092
$ws.Range("J4").Value = 9.961257513846558
$ws.Range("K4").Value = 15.21898526718314
$ws.Range("N4").Value = 17.48621628648429
$ws.Range("O4").Value = 25.53586375545934
$ws.Range("C5").Value = 4.605160051666384
$ws.Range("D5").Value = 8.885640407180601
$ws.Range("E5").Value = 13.48346743079205
$ws.Range("F5").Value = 33.78857725309757
$ws.Range("G5").Value = 3.659065889931818
$ws.Range("J5").Value = 9.964135325371732
$ws.Range("K5").Value = 15.10266227619425
$ws.Range("N5").Value = 17.49842728538048
$ws.Range("O5").Value = 25.55371630187722
$ws.Range("C6").Value = 4.598119124867017
$ws.Range("D6").Value = 8.883885542932203
$ws.Range("E6").Value = 13.48147738528613
$ws.Range("F6").Value = 33.78989680966048
$ws.Range("G6").Value = 3.659179174483187
$ws.Range("J6").Value = 9.96463859193196
$ws.Range("K6").Value = 15.08329068022673
$ws.Range("N6").Value = 17.50048364695192
$ws.Range("O6").Value = 25.55678465619272
$ws.Range("C7").Value = 4.646900739729682
$ws.Range("D7").Value = 8.896182393074783
$ws.Range("E7").Value = 13.49553691804038
$ws.Range("F7").Value = 33.78160296674066
$ws.Range("G7").Value = 3.658400068260646
$ws.Range("J7").Value = 9.961294621677668
$ws.Range("K7").Value = 15.2174203815924
$ws.Range("N7").Value = 17.48637904250103
$ws.Range("O7").Value = 25.53609754791918
$ws.Range("C8").Value = 4.859391310544996
$ws.Range("D8").Value = 8.953265938064835
$ws.Range("E8").Value = 13.56361376713552
$ws.Range("F8").Value = 33.76644795555642
$ws.Range("G8").Value = 3.655140042543036
$ws.Range("J8").Value = 9.950231486607937
$ws.Range("K8").Value = 15.79990447141543
$ws.Range("N8").Value = 17.42830971134472
$ws.Range("O8").Value = 25.45996671050744
$ws.Range("C9").Value = 5.261443841907582
$ws.Range("D9").Value = 9.074794738222765
$ws.Range("E9").Value = 13.71814083728012
$ws.Range("F9").Value = 33.81397547757095
$ws.Range("G9").Value = 3.649382358657098
$ws.Range("J9").Value = 9.941801635677921
$ws.Range("K9").Value = 16.89762816754396
$ws.Range("N9").Value = 17.32937612123297
$ws.Range("O9").Value = 25.36555917918594
$ws.Range("C10").Value = 5.544277735999106
$ws.Range("D10").Value = 9.16940594799334
$ws.Range("E10").Value = 13.84383031914718
$ws.Range("F10").Value = 33.89600353946667
$ws.Range("G10").Value = 3.645535130340999
$ws.Range("J10").Value = 9.943680122281066
$ws.Range("K10").Value = 17.66846693282252
$ws.Range("N10").Value = 17.26575110314621
$ws.Range("O10").Value = 25.32998575286488
$ws.Range("C11").Value = 5.669572453470605
$ws.Range("D11").Value = 9.213485429772643
$ws.Range("E11").Value = 13.90349355057506
$ws.Range("F11").Value = 33.94351544703883
$ws.Range("G11").Value = 3.643867162489959
$ws.Range("J11").Value = 9.946281745470424
$ws.Range("K11").Value = 18.00983743008923
$ws.Range("N11").Value = 17.23876462660627
$ws.Range("O11").Value = 25.32119615609275
$ws.Range("C12").Value = 5.716487738534877
$ws.Range("D12").Value = 9.230316437176549
$ws.Range("E12").Value = 13.92642979119757
$ws.Range("F12").Value = 33.96296719045193
$ws.Range("G12").Value = 3.643247290664433
$ws.Range("J12").Value = 9.947517402451131
$ws.Range("K12").Value = 18.13765571290498
$ws.Range("N12").Value = 17.22882626086918
$ws.Range("O12").Value = 25.31893434137613
$ws.Range("C13").Value = 5.706408002415883
$ws.Range("D13").Value = 9.226685563714533
$ws.Range("E13").Value = 13.92147504002124
$ws.Range("F13").Value = 33.95871309526679
$ws.Range("G13").Value = 3.643380269496527
$ws.Range("J13").Value = 9.947240154927208
$ws.Range("K13").Value = 18.11019402886862
$ws.Range("N13").Value = 17.23095418466943
$ws.Range("O13").Value = 25.31937397749624
$ws.Range("C14").Value = 5.673443095557225
$ws.Range("D14").Value = 9.214867400944376
$ws.Range("E14").Value = 13.90537373827806
$ws.Range("F14").Value = 33.94508653558735
$ws.Range("G14").Value = 3.643815930102368
$ws.Range("J14").Value = 9.946378388538074
$ws.Range("K14").Value = 18.02038282630279
$ws.Range("N14").Value = 17.23794136624823
$ws.Range("O14").Value = 25.32098868116015
$ws.Range("C15").Value = 5.6531806744748
$ws.Range("D15").Value = 9.207646228849438
$ws.Range("E15").Value = 13.89555546410431
$ws.Range("F15").Value = 33.9369297922435
$ws.Range("G15").Value = 3.64408431333595
$ws.Range("J15").Value = 9.945883126829196
$ws.Range("K15").Value = 17.965178626299
$ws.Range("N15").Value = 17.24225777226573
$ws.Range("O15").Value = 25.32211673033925
$ws.Range("C16").Value = 5.536017537314384
$ws.Range("D16").Value = 9.166545289142162
$ws.Range("E16").Value = 13.83997999637788
$ws.Range("F16").Value = 33.89310321855046
$ws.Range("G16").Value = 3.645645783792741
$ws.Range("J16").Value = 9.943545194991508
$ws.Range("K16").Value = 17.64596051863768
$ws.Range("N16").Value = 17.26755406092509
$ws.Range("O16").Value = 25.33070925614147
$ws.Range("C17").Value = 5.46324413466133
$ws.Range("D17").Value = 9.14159022970056
$ws.Range("E17").Value = 13.80651297770976
$ws.Range("F17").Value = 33.86882484225671
$ws.Range("G17").Value = 3.646624692772856
$ws.Range("J17").Value = 9.942557930962707
$ws.Range("K17").Value = 17.44766571769579
$ws.Range("N17").Value = 17.28357330337544
$ws.Range("O17").Value = 25.33787673794365
$ws.Range("C18").Value = 5.421072245967527
$ws.Range("D18").Value = 9.127335309927746
$ws.Range("E18").Value = 13.78749857733761
$ws.Range("F18").Value = 33.85582084326958
$ws.Range("G18").Value = 3.647195471907152
$ws.Range("J18").Value = 9.942154580859212
$ws.Range("K18").Value = 17.33274374743027
$ws.Range("N18").Value = 17.2929713695703
$ws.Range("O18").Value = 25.34269497718569
$ws.Range("C19").Value = 5.406741059313327
$ws.Range("D19").Value = 9.122526086063106
$ws.Range("E19").Value = 13.78110141120583
$ws.Range("F19").Value = 33.85158300240609
$ws.Range("G19").Value = 3.647390058590991
$ws.Range("J19").Value = 9.942046286758123
$ws.Range("K19").Value = 17.2936877685718
$ws.Range("N19").Value = 17.29618505122938
$ws.Range("O19").Value = 25.34444571695126
$ws.Range("C20").Value = 5.471023898748402
$ws.Range("D20").Value = 9.144236612373884
$ws.Range("E20").Value = 13.8100513862082
$ws.Range("F20").Value = 33.87130997085976
$ws.Range("G20").Value = 3.646519685987219
$ws.Range("J20").Value = 9.94264600697128
$ws.Range("K20").Value = 17.46886525155689
$ws.Range("N20").Value = 17.28184896558774
$ws.Range("O20").Value = 25.33704171773154
$ws.Range("C21").Value = 5.683140448501101
$ws.Range("D21").Value = 9.218334992703026
$ws.Range("E21").Value = 13.91009388837724
$ws.Range("F21").Value = 33.94904942163754
$ws.Range("G21").Value = 3.643687647588081
$ws.Range("J21").Value = 9.946624718963335
$ws.Range("K21").Value = 18.046802807999
$ws.Range("N21").Value = 17.23588144448263
$ws.Range("O21").Value = 25.32048543155279
$ws.Range("C22").Value = 5.818657158542799
$ws.Range("D22").Value = 9.267568111895407
$ws.Range("E22").Value = 13.97747000549631
$ws.Range("F22").Value = 34.00836247680807
$ws.Range("G22").Value = 3.641905214247969
$ws.Range("J22").Value = 9.950684573105724
$ws.Range("K22").Value = 18.4160132266349
$ws.Range("N22").Value = 17.20747564122662
$ws.Range("O22").Value = 25.31588309451735
$ws.Range("C23").Value = 5.746628350837246
$ws.Range("D23").Value = 9.241221296698567
$ws.Range("E23").Value = 13.941332715691
$ws.Range("F23").Value = 33.97593026741616
$ws.Range("G23").Value = 3.642850288219403
$ws.Range("J23").Value = 9.94838448866494
$ws.Range("K23").Value = 18.21977223429096
$ws.Range("N23").Value = 17.22248678287908
$ws.Range("O23").Value = 25.317769516586
$ws.Range("C24").Value = 5.467507702578338
$ws.Range("D24").Value = 9.143039895174084
$ws.Range("E24").Value = 13.80845096595482
$ws.Range("F24").Value = 33.87018347214912
$ws.Range("G24").Value = 3.64656713465732
$ws.Range("J24").Value = 9.942605676156143
$ws.Range("K24").Value = 17.45928379785379
$ws.Range("N24").Value = 17.28262795180566
$ws.Range("O24").Value = 25.33741705771497
$ws.Range("C25").Value = 5.154643385176754
$ws.Range("D25").Value = 9.040943942996281
$ws.Range("E25").Value = 13.67415111433777
$ws.Range("F25").Value = 33.79283911167211
$ws.Range("G25").Value = 3.650872406821493
$ws.Range("J25").Value = 9.942663040418962
$ws.Range("K25").Value = 16.60635735598012
$ws.Range("N25").Value = 17.3545457183702
$ws.Range("O25").Value = 25.38518664982512
